$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 147, shifting existing rows 147-221 down to 148-222.
$ws.Rows("147:147").Insert()

# Populate the newly inserted row 147 with the new record's data.
$ws.Range("A147").Value = 5
$ws.Range("B147").Value = "Macroferia Regional de Talca"
$ws.Range("C147").Value = "Maule"
$ws.Range("D147").Value = 44466
$ws.Range("E147").Value = 7
$ws.Range("F147").Value = 100112043
$ws.Range("G147").Value = "Pepino ensalada"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 300
$ws.Range("K147").Value = 16000
$ws.Range("L147").Value = 16000
$ws.Range("M147").Value = 16000
$ws.Range("N147").Value = "`$/caja 60 unidades"
$ws.Range("O147").Value = "Región de Arica y Parinacota"
$ws.Range("P147").Value = 267
$ws.Range("Q147").Value = 60
$ws.Range("R147").Value = "Hortaliza"
